# Daily attendance processing - 2025-12-26 17:55:43
#
# The "Recorded By" column (G) holds a comma-separated list of the
# users/processes that touched each attendance session (e.g.
# "dnasr281@gmail.com, System"). This pass normalizes the ordering by
# swapping the last two entries in that list for every row that has
# more than one recorder, leaving single-recorder rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        $n = $parts.Count

        if ($n -ge 2) {
            $last = $parts[$n - 1]
            $parts[$n - 1] = $parts[$n - 2]
            $parts[$n - 2] = $last

            $cell.Value = $parts -join ", "
        }
    }
}
